$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 12) to the SO_Test sheet, mirroring the
# unstyled look of rows 8-11 (no explicit cell style, default format).
#
# Columns A and D hold values that look like a date ("01/12/2024") and a
# long numeric phone string ("9157994875"); assigning them straight would
# make Excel auto-convert them into a date serial number / plain number.
# Prefixing with a leading apostrophe forces them to be stored as literal
# text, matching the source data (t="inlineStr"/shared-string "a"-style
# entries in the target workbook).
$ws.Cells.Item(12, 1).Value = "'01/12/2024"
$ws.Cells.Item(12, 2).Value = "SO240112006"
$ws.Cells.Item(12, 3).Value = "a"
$ws.Cells.Item(12, 4).Value = "'9157994875"
$ws.Cells.Item(12, 5).Value = "a"
$ws.Cells.Item(12, 6).Value = "a"
$ws.Cells.Item(12, 7).Value = 2
$ws.Cells.Item(12, 8).Value = 42
$ws.Cells.Item(12, 9).Value = "AEC"
$ws.Cells.Item(12, 10).Value = "OTHER"
$ws.Cells.Item(12, 11).Value = "a"

# Writing the values above causes Excel to stamp the new cells with the
# worksheet's per-column style (and/or a quote-prefix style for the two
# text-forced cells). Rows 8-11 in the source workbook carry no explicit
# style at all, so copy the (lack of) formatting from row 11 back onto
# row 12 to match - this only touches formatting, not the values just set.
$ws.Range("A11:K11").Copy()
$ws.Range("A12:K12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
